# Update the "Data information:" description text on the Metadata sheet.
# The description now clarifies that suppressed *regions* (in addition to
# suppressed categories) are excluded from the percentage calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B10").Value = "Official annual DBT estimates of inward foreign direct investment (FDI) activity, reporting the percentage of FDI projects from new investments, and excisiting investment. Percentages are calculated from published project counts only and exclude suppressed categories and regions. Sub-national figures exclude multi-site FDI projects. Some regional values are suppressed to protect confidentiality."

$ws.Range("B11").Select()
